$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsFoF = $wb.Worksheets.Item("FoFEtiL")

# Update the label in the FoFEtiL sheet to clarify the unit is dimensionless,
# and make the header wrap / grow the row to fit the longer text.
$wsFoF.Range("B1").Value = "Fraction of Forestry Expenses (dimensionless)"
$wsFoF.Range("B1").WrapText = $true
$wsFoF.Rows.Item(1).RowHeight = 28.5

# Reflect the new active cell selection on the FoFEtiL sheet without changing
# which sheet tab is active (the "About" sheet stays the selected tab).
$wsFoF.Activate()
[void]$wsFoF.Range("B1").Select()
$wsAbout.Activate()
